$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Ntn1"
$ws.Range("C2").Value = "Unc5d"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 2.8328235
$ws.Range("H2").Value = 5.665647
$ws.Range("I2").Value = 0.1613214142302211
$ws.Range("J2").Value = 0.1283172186594826
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.5
$ws.Range("M2").Value = 0.012541
$ws.Range("N2").Value = 0.025082
$ws.Range("O2").Value = 0.07625098802213169
$ws.Range("P2").Value = 0.07625098802213169
$ws.Range("Q2").Value = 0.0355264395135
$ws.Range("R2").Value = 0.142105758054
$ws.Range("S2").Value = 0.01230091722418193
$ws.Range("T2").Value = 0.009784314703037459

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Ntn1"
$ws.Range("C3").Value = "Unc5d"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 2.8328235
$ws.Range("H3").Value = 5.665647
$ws.Range("I3").Value = 0.1613214142302211
$ws.Range("J3").Value = 0.1283172186594826
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.151929
$ws.Range("N3").Value = 0.303858
$ws.Range("O3").Value = 0.9237490119778683
$ws.Range("P3").Value = 0.9237490119778683
$ws.Range("Q3").Value = 0.4303880415315
$ws.Range("R3").Value = 1.721552166126
$ws.Range("S3").Value = 0.1490204970060392
$ws.Range("T3").Value = 0.1185329039564451

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Ntn1"
$ws.Range("C4").Value = "Unc5d"
$ws.Range("D4").Value = "ECs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 8.502875000000001
$ws.Range("H4").Value = 25.508625
$ws.Range("I4").Value = 0.4842150667074004
$ws.Range("J4").Value = 0.5777267471531042
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.5
$ws.Range("M4").Value = 0.012541
$ws.Range("N4").Value = 0.025082
$ws.Range("O4").Value = 0.07625098802213169
$ws.Range("P4").Value = 0.07625098802213169
$ws.Range("Q4").Value = 0.106634555375
$ws.Range("R4").Value = 0.6398073322500001
$ws.Range("S4").Value = 0.03692187725164169
$ws.Range("T4").Value = 0.04405223527723646

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Ntn1"
$ws.Range("C5").Value = "Unc5d"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 8.502875000000001
$ws.Range("H5").Value = 25.508625
$ws.Range("I5").Value = 0.4842150667074004
$ws.Range("J5").Value = 0.5777267471531042
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.151929
$ws.Range("N5").Value = 0.303858
$ws.Range("O5").Value = 0.9237490119778683
$ws.Range("P5").Value = 0.9237490119778683
$ws.Range("Q5").Value = 1.291833295875
$ws.Range("R5").Value = 7.750999775250001
$ws.Range("S5").Value = 0.4472931894557587
$ws.Range("T5").Value = 0.5336745118758678

# Row 6
$ws.Range("A6").Value = "Inflammatory-Mac"
$ws.Range("B6").Value = "Ntn1"
$ws.Range("C6").Value = "Unc5d"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.1821286666666667
$ws.Range("H6").Value = 0.546386
$ws.Range("I6").Value = 0.01037172068027931
$ws.Range("J6").Value = 0.01237470880809906
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.5
$ws.Range("M6").Value = 0.012541
$ws.Range("N6").Value = 0.025082
$ws.Range("O6").Value = 0.07625098802213169
$ws.Range("P6").Value = 0.07625098802213169
$ws.Range("Q6").Value = 0.002284075608666667
$ws.Range("R6").Value = 0.013704453652
$ws.Range("S6").Value = 0.0007908539493608729
$ws.Range("T6").Value = 0.0009435837731037294

# Row 7
$ws.Range("A7").Value = "Inflammatory-Mac"
$ws.Range("B7").Value = "Ntn1"
$ws.Range("C7").Value = "Unc5d"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.1821286666666667
$ws.Range("H7").Value = 0.546386
$ws.Range("I7").Value = 0.01037172068027931
$ws.Range("J7").Value = 0.01237470880809906
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.151929
$ws.Range("N7").Value = 0.303858
$ws.Range("O7").Value = 0.9237490119778683
$ws.Range("P7").Value = 0.9237490119778683
$ws.Range("Q7").Value = 0.02767062619800001
$ws.Range("R7").Value = 0.166023757188
$ws.Range("S7").Value = 0.009580866730918434
$ws.Range("T7").Value = 0.01143112503499534

# Row 8
$ws.Range("A8").Value = "MuSCs"
$ws.Range("B8").Value = "Ntn1"
$ws.Range("C8").Value = "Unc5d"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 5.6940955
$ws.Range("H8").Value = 11.388191
$ws.Range("I8").Value = 0.3242628913597822
$ws.Range("J8").Value = 0.257923057098854
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.5
$ws.Range("M8").Value = 0.012541
$ws.Range("N8").Value = 0.025082
$ws.Range("O8").Value = 0.07625098802213169
$ws.Range("P8").Value = 0.07625098802213169
$ws.Range("Q8").Value = 0.07140965166549999
$ws.Range("R8").Value = 0.285638606662
$ws.Range("S8").Value = 0.02472536584509654
$ws.Range("T8").Value = 0.01966688793747631

# Row 9
$ws.Range("A9").Value = "MuSCs"
$ws.Range("B9").Value = "Ntn1"
$ws.Range("C9").Value = "Unc5d"
$ws.Range("D9").Value = "MuSCs"
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 5.6940955
$ws.Range("H9").Value = 11.388191
$ws.Range("I9").Value = 0.3242628913597822
$ws.Range("J9").Value = 0.257923057098854
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.151929
$ws.Range("N9").Value = 0.303858
$ws.Range("O9").Value = 0.9237490119778683
$ws.Range("P9").Value = 0.9237490119778683
$ws.Range("Q9").Value = 0.8650982352195
$ws.Range("R9").Value = 3.460392940878
$ws.Range("S9").Value = 0.2995375255146857
$ws.Range("T9").Value = 0.2382561691613777

# Row 10
$ws.Range("A10").Value = "Neutrophils"
$ws.Range("B10").Value = "Ntn1"
$ws.Range("C10").Value = "Unc5d"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 0.6666666666666666
$ws.Range("G10").Value = 0.165052
$ws.Range("H10").Value = 0.495156
$ws.Range("I10").Value = 0.009399252040067608
$ws.Range("J10").Value = 0.01121443689000652
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.5
$ws.Range("M10").Value = 0.012541
$ws.Range("N10").Value = 0.025082
$ws.Range("O10").Value = 0.07625098802213169
$ws.Range("P10").Value = 0.07625098802213169
$ws.Range("Q10").Value = 0.002069917132
$ws.Range("R10").Value = 0.012419502792
$ws.Range("S10").Value = 0.000716702254724192
$ws.Range("T10").Value = 0.0008551118929748386

# Row 11
$ws.Range("A11").Value = "Neutrophils"
$ws.Range("B11").Value = "Ntn1"
$ws.Range("C11").Value = "Unc5d"
$ws.Range("D11").Value = "MuSCs"
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 0.6666666666666666
$ws.Range("G11").Value = 0.165052
$ws.Range("H11").Value = 0.495156
$ws.Range("I11").Value = 0.009399252040067608
$ws.Range("J11").Value = 0.01121443689000652
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 0.151929
$ws.Range("N11").Value = 0.303858
$ws.Range("O11").Value = 0.9237490119778683
$ws.Range("P11").Value = 0.9237490119778683
$ws.Range("Q11").Value = 0.025076185308
$ws.Range("R11").Value = 0.150457111848
$ws.Range("S11").Value = 0.008682549785343415
$ws.Range("T11").Value = 0.01035932499703168

# Row 12
$ws.Range("A12").Value = "Resolving-Mac"
$ws.Range("B12").Value = "Ntn1"
$ws.Range("C12").Value = "Unc5d"
$ws.Range("D12").Value = "ECs"
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 0.6666666666666666
$ws.Range("G12").Value = 0.183146
$ws.Range("H12").Value = 0.549438
$ws.Range("I12").Value = 0.01042965498224936
$ws.Range("J12").Value = 0.01244383139045351
$ws.Range("K12").Value = 1
$ws.Range("L12").Value = 0.5
$ws.Range("M12").Value = 0.012541
$ws.Range("N12").Value = 0.025082
$ws.Range("O12").Value = 0.07625098802213169
$ws.Range("P12").Value = 0.07625098802213169
$ws.Range("Q12").Value = 0.002296833986
$ws.Range("R12").Value = 0.013781003916
$ws.Range("S12").Value = 0.0007952714971264623
$ws.Range("T12").Value = 0.0009488544383028973

# Row 13
$ws.Range("A13").Value = "Resolving-Mac"
$ws.Range("B13").Value = "Ntn1"
$ws.Range("C13").Value = "Unc5d"
$ws.Range("D13").Value = "MuSCs"
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 0.6666666666666666
$ws.Range("G13").Value = 0.183146
$ws.Range("H13").Value = 0.549438
$ws.Range("I13").Value = 0.01042965498224936
$ws.Range("J13").Value = 0.01244383139045351
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.151929
$ws.Range("N13").Value = 0.303858
$ws.Range("O13").Value = 0.9237490119778683
$ws.Range("P13").Value = 0.9237490119778683
$ws.Range("Q13").Value = 0.027825188634
$ws.Range("R13").Value = 0.166951131804
$ws.Range("S13").Value = 0.009634383485122902
$ws.Range("T13").Value = 0.01149497695215062
